$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42: Eye of the Beholder
$ws.Range("H42").Value = 264.9
$ws.Range("I42").Value = 20.833334
$ws.Range("K42").Value = 62.500002
$ws.Range("M42").Value = 167.499998

# Row 69: Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 7071.4287
$ws.Range("I69").Value = 3500
$ws.Range("J69").Value = 7666.6665
$ws.Range("K69").Value = 10500
$ws.Range("L69").Value = 22999.9995
$ws.Range("M69").Value = -9626
$ws.Range("N69").Value = -24747.9995

# Row 72: Surgical Substitution (L)
$ws.Range("H72").Value = 7071.4287
$ws.Range("I72").Value = 3500
$ws.Range("J72").Value = 7666.6665
$ws.Range("K72").Value = 31500
$ws.Range("L72").Value = 68999.9985
$ws.Range("M72").Value = -27132
$ws.Range("N72").Value = -77735.9985

# Row 82: Rolling on Initiative
$ws.Range("H82").Value = 7835.125

# Row 85: Darkly Dreaming Dexterity (L)
$ws.Range("H85").Value = 7835.125

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 3818.6287
$ws.Range("J112").Value = 4185.1
$ws.Range("L112").Value = 12555.3
$ws.Range("N112").Value = -14771.3

# Row 131: Mindful Study
$ws.Range("H131").Value = 3732.5
$ws.Range("I131").Value = 841.4545000000001
$ws.Range("J131").Value = 14333
$ws.Range("K131").Value = 2524.3635
$ws.Range("L131").Value = 42999
$ws.Range("M131").Value = 2515.6365
$ws.Range("N131").Value = -53079

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 18820.076
$ws.Range("I132").Value = 5440
$ws.Range("K132").Value = 16320
$ws.Range("M132").Value = -13790

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 11496880
$ws.Range("I137").Value = 1933.2142
$ws.Range("J137").Value = 22225498
$ws.Range("K137").Value = 5799.642599999999
$ws.Range("L137").Value = 66676494
$ws.Range("M137").Value = -3249.642599999999
$ws.Range("N137").Value = -66681594

# Row 138: All-night Crafting
$ws.Range("H138").Value = 6431.3096
$ws.Range("J138").Value = 7084.4053
$ws.Range("L138").Value = 21253.2159
$ws.Range("N138").Value = -31533.2159

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2709
$ws.Range("I45").Value = 2001.6666
$ws.Range("J45").Value = 2901.9092
$ws.Range("K45").Value = 2001.6666
$ws.Range("L45").Value = 2901.9092
$ws.Range("M45").Value = -1624.6666
$ws.Range("N45").Value = -3655.9092

# Row 97: Ore for Me
$ws.Range("H97").Value = 531.9286
$ws.Range("I97").Value = 304.9
$ws.Range("K97").Value = 304.9
$ws.Range("M97").Value = 191.1

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience
$ws.Range("H107").Value = 2511.2354
$ws.Range("I107").Value = 2835.8333
$ws.Range("K107").Value = 2835.8333
$ws.Range("M107").Value = -915.8332999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 93: Reeling for Rods
$ws.Range("H93").Value = 142877440
$ws.Range("J93").Value = 333369660
$ws.Range("L93").Value = 333369660
$ws.Range("N93").Value = -333373404

# Row 96: Composition
$ws.Range("H96").Value = 14999.6
$ws.Range("J96").Value = 14999.6
$ws.Range("L96").Value = 14999.6
$ws.Range("N96").Value = -20491.6

# Row 99: O Pine
$ws.Range("H99").Value = 8057
$ws.Range("I99").Value = 3402.4
$ws.Range("K99").Value = 3402.4
$ws.Range("M99").Value = -1904.4

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1894758.5
$ws.Range("I105").Value = 2841487.5
$ws.Range("J105").Value = 1300.5
$ws.Range("K105").Value = 2841487.5
$ws.Range("L105").Value = 1300.5
$ws.Range("M105").Value = -2839740.5
$ws.Range("N105").Value = -4794.5

# Row 107: Built to Last
$ws.Range("H107").Value = 957712.25
$ws.Range("I107").Value = 2273196
$ws.Range("K107").Value = 2273196
$ws.Range("M107").Value = -2271276

# Row 126: A Better Conductor
$ws.Range("H126").Value = 8057
$ws.Range("I126").Value = 3402.4
$ws.Range("K126").Value = 10207.2
$ws.Range("M126").Value = -7737.200000000001

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 83342620
$ws.Range("J132").Value = 39999.668
$ws.Range("L132").Value = 119999.004
$ws.Range("N132").Value = -125059.004

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2886.3635
$ws.Range("I134").Value = 2605.8
$ws.Range("K134").Value = 7817.400000000001
$ws.Range("M134").Value = -5282.400000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 56: Culture Club
$ws.Range("H56").Value = 55563690
$ws.Range("I56").Value = 55563690
$ws.Range("K56").Value = 55563690
$ws.Range("M56").Value = -55563160

# Row 57: The Egg Files
$ws.Range("H57").Value = 5265
$ws.Range("J57").Value = 12500
$ws.Range("L57").Value = 37500
$ws.Range("N57").Value = -38618

# Row 87: Soup That Eats Like a Knight
$ws.Range("H87").Value = 12052.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 12052.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 36157.5
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -38653.5

# Row 90: Like Ma Used to Make (L)
$ws.Range("H90").Value = 12052.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 12052.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 108472.5
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -120952.5

# Row 107: Slippery Service
$ws.Range("H107").Value = 387.44446
$ws.Range("J107").Value = 470.5
$ws.Range("L107").Value = 1411.5
$ws.Range("N107").Value = -5251.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 1513805.8
$ws.Range("I80").Value = 1849488.8
$ws.Range("J80").Value = 3232
$ws.Range("K80").Value = 1849488.8
$ws.Range("L80").Value = 3232
$ws.Range("M80").Value = -1848490.8
$ws.Range("N80").Value = -5228

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1513805.8
$ws.Range("I83").Value = 1849488.8
$ws.Range("J83").Value = 3232
$ws.Range("K83").Value = 9247444
$ws.Range("L83").Value = 16160
$ws.Range("M83").Value = -9242452
$ws.Range("N83").Value = -26144

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 700.7059
$ws.Range("J97").Value = 812.6
$ws.Range("L97").Value = 812.6
$ws.Range("N97").Value = -1804.6

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 337798.9
$ws.Range("I122").Value = 669101.7
$ws.Range("K122").Value = 2007305.1
$ws.Range("M122").Value = -2004855.1

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 4505.8423
$ws.Range("J126").Value = 6143.727
$ws.Range("L126").Value = 18431.181
$ws.Range("N126").Value = -23371.181

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3594.2856
$ws.Range("I132").Value = 3120.0908
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 9360.2724
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -6830.2724
$ws.Range("N132").Value = -21059

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 7816612
$ws.Range("I40").Value = 9618599
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 9618599
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -9618463
$ws.Range("N40").Value = -8272

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 1564215.2
$ws.Range("J82").Value = 2457.889
$ws.Range("L82").Value = 2457.889
$ws.Range("N82").Value = -3179.889

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 1564215.2
$ws.Range("J85").Value = 2457.889
$ws.Range("L85").Value = 2457.889
$ws.Range("N85").Value = -4953.889

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1496.0454
$ws.Range("I93").Value = 1372.7858
$ws.Range("J93").Value = 1711.75
$ws.Range("K93").Value = 1372.7858
$ws.Range("L93").Value = 1711.75
$ws.Range("M93").Value = -124.7858000000001
$ws.Range("N93").Value = -4207.75

# Row 98: Try Tricorne Again
$ws.Range("H98").Value = 49987
$ws.Range("J98").Value = 49987
$ws.Range("L98").Value = 49987
$ws.Range("N98").Value = -55977

# Row 122: Hell on Leather
$ws.Range("H122").Value = 8090.6895
$ws.Range("I122").Value = 3870.0715
$ws.Range("J122").Value = 12029.934
$ws.Range("K122").Value = 11610.2145
$ws.Range("L122").Value = 36089.802
$ws.Range("M122").Value = -9160.2145
$ws.Range("N122").Value = -40989.802

# Row 141: Just Generally Freezing
$ws.Range("H141").Value = 97950
$ws.Range("J141").Value = 97950
$ws.Range("L141").Value = 97950
$ws.Range("N141").Value = -108310

$ws = $wb.Worksheets.Item("WVR")
# Row 61: Bundle Up, It's Odd out There
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 69: Fashion Patrol
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72: Dress Code Violation (L)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 767
$ws.Range("I107").Value = 492.125
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1476.375
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = 443.625
$ws.Range("N107").Value = -8340

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 4143.95
$ws.Range("I122").Value = 2536.923
$ws.Range("K122").Value = 7610.768999999999
$ws.Range("M122").Value = -5160.768999999999
